# This script applies a "base update" to the Montenegro Prva Liga sheet.
# The underlying change re-associates several fixture rows with the
# correct match data: rows 29/30, 75/76, 162/163 swap their payload
# (everything except the running id in col A, the Div in col C and the
# Date in col D), while rows 153/154/155 rotate their payload in a
# three-way cycle (153<-154<-155<-153). The id numbers in column B make
# the intended final arrangement unambiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Montenegro Prva Liga")

# Target values keyed by row number, then by 1-based column index.
# Column layout: B=2 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15
#                P=16 Q=17 R=18 S=19 T=20 U=21 V=22 W=23 X=24 Y=25
#                Z=26 AA=27 AB=28
$data = @{
  29 = @{
    2 = 6815315
    5 = 'FK Decic Tuzi'
    6 = 'FK Rudar Pljevlja'
    7 = 0
    8 = 1
    9 = 'A'
    10 = 1.615
    11 = 3.5
    12 = 4.75
    13 = 1.4
    14 = 3.8
    15 = 6.5
    16 = -1.25
    17 = 2
    18 = 1.8
    19 = 2.5
    20 = 1.95
    21 = 1.85
    22 = -1
    23 = -1
    24 = 5.5
    25 = -1
    26 = 0.8
    27 = -1
    28 = 0.8500000000000001
  }
  30 = @{
    2 = 6815312
    5 = 'Buducnost Podgorica'
    6 = 'FK Arsenal'
    7 = 2
    8 = 0
    9 = 'H'
    10 = 1.444
    11 = 4
    12 = 6
    13 = 1.4
    14 = 4
    15 = 6.5
    16 = -1.25
    17 = 1.95
    18 = 1.85
    19 = 2.5
    20 = 1.775
    21 = 1.925
    22 = 0.3999999999999999
    23 = -1
    24 = -1
    25 = 0.95
    26 = -1
    27 = -1
    28 = 0.925
  }
  75 = @{
    2 = 6815359
    5 = 'Buducnost Podgorica'
    6 = 'FK Jedinstvo Bijelo Polje'
    7 = 3
    8 = 2
    9 = 'H'
    10 = 1.333
    11 = 4.333
    12 = 7.5
    13 = 1.333
    14 = 4.333
    15 = 8
    16 = -1.5
    17 = 1.875
    18 = 1.925
    19 = 2.75
    20 = 1.8
    21 = 2
    22 = 0.333
    23 = -1
    24 = -1
    25 = -1
    26 = 0.925
    27 = 0.8
    28 = -1
  }
  76 = @{
    2 = 6815358
    5 = 'OFK Petrovac'
    6 = 'FK Arsenal'
    7 = 1
    8 = 1
    9 = 'D'
    10 = 2.1
    11 = 3.1
    12 = 3.2
    13 = 1.75
    14 = 3.3
    15 = 4.2
    16 = -0.5
    17 = 1.8
    18 = 2
    19 = 2.25
    20 = 1.95
    21 = 1.85
    22 = -1
    23 = 2.3
    24 = -1
    25 = -1
    26 = 1
    27 = -0.5
    28 = 0.425
  }
  153 = @{
    2 = 6816282
    5 = 'Buducnost Podgorica'
    6 = 'Sutjeska Niksic'
    7 = 2
    8 = 2
    9 = 'D'
    10 = 2.6
    11 = 3
    12 = 2.5
    13 = 1.5
    14 = 3.8
    15 = 5.25
    16 = -1
    17 = 1.95
    18 = 1.85
    19 = 2.5
    20 = 1.75
    21 = 1.95
    22 = -1
    23 = 2.8
    24 = -1
    25 = -1
    26 = 0.8500000000000001
    27 = 0.75
    28 = -1
  }
  154 = @{
    2 = 6816283
    5 = 'FK Jedinstvo Bijelo Polje'
    6 = 'FK Arsenal'
    7 = 1
    8 = 1
    9 = 'D'
    10 = 2.15
    11 = 3
    12 = 3.2
    13 = 1.95
    14 = 3.1
    15 = 3.5
    16 = -0.25
    17 = 1.725
    18 = 1.975
    19 = 2.25
    20 = 1.825
    21 = 1.975
    22 = -1
    23 = 2.1
    24 = -1
    25 = -0.5
    26 = 0.4875
    27 = -0.5
    28 = 0.4875
  }
  155 = @{
    2 = 6816281
    5 = 'OFK Petrovac'
    6 = 'FK Jezero'
    7 = 2
    8 = 1
    9 = 'H'
    10 = 2.25
    11 = 3
    12 = 3
    13 = 2.3
    14 = 2.875
    15 = 2.9
    16 = -0.25
    17 = 2.05
    18 = 1.75
    19 = 1.75
    20 = 1.875
    21 = 1.925
    22 = 1.3
    23 = -1
    24 = -1
    25 = 1.05
    26 = -1
    27 = 0.875
    28 = -1
  }
  162 = @{
    2 = 6817582
    5 = 'OFK Mladost DG'
    6 = 'Sutjeska Niksic'
    7 = 3
    8 = 2
    9 = 'H'
    10 = 4.2
    11 = 3.3
    12 = 1.75
    13 = 4.2
    14 = 3.3
    15 = 1.75
    16 = 0.5
    17 = 2
    18 = 1.8
    19 = 2.25
    20 = 1.8
    21 = 2
    22 = 3.2
    23 = -1
    24 = -1
    25 = 1
    26 = -1
    27 = 0.8
    28 = -1
  }
  163 = @{
    2 = 6817583
    5 = 'OFK Petrovac'
    6 = 'FK Arsenal'
    7 = 2
    8 = 2
    9 = 'D'
    10 = 1.75
    11 = 3.3
    12 = 4.2
    13 = 2.05
    14 = 3.3
    15 = 3.1
    16 = -0.25
    17 = 1.875
    18 = 1.925
    19 = 2.5
    20 = 1.85
    21 = 1.95
    22 = -1
    23 = 2.3
    24 = -1
    25 = -0.5
    26 = 0.4625
    27 = 0.8500000000000001
    28 = -1
  }
}

foreach ($rowNum in $data.Keys) {
  $cols = $data[$rowNum]
  foreach ($colIdx in $cols.Keys) {
    $ws.Cells.Item($rowNum, $colIdx).Value = $cols[$colIdx]
  }
}
